$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 113, shifting existing rows 113.. down by one.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new record's data.
$ws.Cells.Item(113, 1).Value = 1
$ws.Cells.Item(113, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(113, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(113, 4).Value = 45068
$ws.Cells.Item(113, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(113, 5).Value = 15
$ws.Cells.Item(113, 6).Value = "Fruta"
$ws.Cells.Item(113, 7).Value = 100106
$ws.Cells.Item(113, 8).Value = "Oleaginosos"
$ws.Cells.Item(113, 9).Value = 100106002
$ws.Cells.Item(113, 10).Value = "Palta"
$ws.Cells.Item(113, 11).Value = "Hass"
$ws.Cells.Item(113, 12).Value = "Segunda"
$ws.Cells.Item(113, 13).Value = 400
$ws.Cells.Item(113, 14).Value = 25000
$ws.Cells.Item(113, 15).Value = 26000
$ws.Cells.Item(113, 16).Value = 25500
$ws.Cells.Item(113, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(113, 18).Value = "Perú"
$ws.Cells.Item(113, 19).Value = 2550
$ws.Cells.Item(113, 20).Value = 10
